# Natmi following Dr Hou advice
# Rebuild the LR-pair table for Reln/Lrp8 so that every Sending-cluster /
# Target-cluster combination across the three clusters (ECs, FAPs, sCs) is
# present (previously only 6 of the 9 combinations were listed), and update
# the numeric NATMI statistics accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [string]$Sending,
        [string]$Ligand,
        [string]$Receptor,
        [string]$Target,
        [double[]]$Values
    )

    $ws.Cells.Item($Row, 1).Value = $Sending
    $ws.Cells.Item($Row, 2).Value = $Ligand
    $ws.Cells.Item($Row, 3).Value = $Receptor
    $ws.Cells.Item($Row, 4).Value = $Target

    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, 5 + $i).Value = $Values[$i]
    }
}

Set-Row 2 "ECs" "Reln" "Lrp8" "ECs" @(2,0.6666666666666666,0.2638773333333334,0.791632,0.0472992197617186,0.04729921976171859,3,1,2.228108666666667,6.684326,0.5687849952918405,0.5687849952918405,0.587947373336889,5.291526360032,0.02690308648947684,0.02690308648947684)
Set-Row 3 "ECs" "Reln" "Lrp8" "FAPs" @(2,0.6666666666666666,0.2638773333333334,0.791632,0.0472992197617186,0.04729921976171859,3,1,0.1124773333333333,0.337432,0.02871288122861097,0.02871288122861097,0.02968021878044444,0.267121969024,0.001358096879224195,0.001358096879224195)
Set-Row 4 "ECs" "Reln" "Lrp8" "sCs" @(2,0.6666666666666666,0.2638773333333334,0.791632,0.0472992197617186,0.04729921976171859,3,1,1.576726666666667,4.73018,0.4025021234795486,0.4025021234795487,0.4160624281955556,3.74456185376,0.01903803639301757,0.01903803639301757)
Set-Row 5 "FAPs" "Reln" "Lrp8" "ECs" @(3,1,1.127914,3.383742,0.2021751981665183,0.2021751981665183,3,1,2.228108666666667,6.684326,0.5687849952918405,0.5687849952918405,2.513114958654667,22.618034627892,0.11499421913727,0.11499421913727)
Set-Row 6 "FAPs" "Reln" "Lrp8" "FAPs" @(3,1,1.127914,3.383742,0.2021751981665183,0.2021751981665183,3,1,0.1124773333333333,0.337432,0.02871288122861097,0.02871288122861097,0.1268647589493333,1.141782830544,0.005805032452326125,0.005805032452326126)
Set-Row 7 "FAPs" "Reln" "Lrp8" "sCs" @(3,1,1.127914,3.383742,0.2021751981665183,0.2021751981665183,3,1,1.576726666666667,4.73018,0.4025021234795486,0.4025021234795487,1.778412081506666,16.00570873356,0.08137594657692215,0.08137594657692215)
Set-Row 8 "sCs" "Reln" "Lrp8" "ECs" @(3,1,4.187102666666667,12.561308,0.7505255820717631,0.750525582071763,3,1,2.228108666666667,6.684326,0.5687849952918405,0.5687849952918405,9.329319739823113,83.96387765840801,0.4268876896650936,0.4268876896650936)
Set-Row 9 "sCs" "Reln" "Lrp8" "FAPs" @(3,1,4.187102666666667,12.561308,0.7505255820717631,0.750525582071763,3,1,0.1124773333333333,0.337432,0.02871288122861097,0.02871288122861097,0.4709541423395555,4.238587281056,0.02154975189706065,0.02154975189706065)
Set-Row 10 "sCs" "Reln" "Lrp8" "sCs" @(3,1,4.187102666666667,12.561308,0.7505255820717631,0.750525582071763,3,1,1.576726666666667,4.73018,0.4025021234795486,0.4025021234795487,6.601916430604445,59.41724787544,0.3020881405096089,0.3020881405096089)
